# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, matching the regenerated data dump.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3508
$ws1.Range("F5").Value = 7013
$ws1.Range("F6").Value = 2873
$ws1.Range("F7").Value = 53
$ws1.Range("F8").Value = 130
$ws1.Range("F13").Value = 10
$ws1.Range("F15").Value = 587
$ws1.Range("F16").Value = 12

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3508
$ws4.Range("F6").Value = 7013
$ws4.Range("F7").Value = 2873
$ws4.Range("F8").Value = 53
$ws4.Range("F9").Value = 130
$ws4.Range("F14").Value = 10
$ws4.Range("F16").Value = 587
$ws4.Range("F17").Value = 12
